$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 179878
$ws.Range("C4").Value = 169835
$ws.Range("C7").Value = 5.58
$ws.Range("C8").Value = 65.06999999999999
